# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: copy the existing header style (bold/border/centered) onto
#     the three new header cells, then set their captions. ---
$ws.Range("A1").Copy()
$ws.Range("AD1").PasteSpecial(-4122)
$ws.Range("AE1").PasteSpecial(-4122)
$ws.Range("AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows: every player row (2-48) gets the team's season record. ---
$wins = 90
$losses = 72
$ties = 0

for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
